# Update the "cryptos" price/volume snapshot (GitHub Actions refresh).
# Price (column D) values are kept as plain text (e.g. "3.010.45") just like
# the source data, so every Price write forces Text format ("@") before
# assigning the string and restores General/Normal style afterwards -
# otherwise Excel would "helpfully" reinterpret strings such as "6.05" or
# "0.999" as numbers and mangle the "x.xxx.xx"-style big-number strings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.609.67"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.05%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.010.45"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.11%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.44"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.18"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.55%  "

$ws.Range("E7").Value = "  +0.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.002.67"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.27%  "

$ws.Range("E9").Value = "  -0.67%  "

$ws.Range("E10").Value = "  -5.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.05"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -8.97%  "

$ws.Range("E12").Value = "  -2.31%  "

$ws.Range("E13").Value = "  -2.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.58"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.505.10"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.809.66"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.80%  "

$ws.Range("E17").Value = "  -2.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.009.77"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.66"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.70%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "472.76"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.67%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.28"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.67%  "

$ws.Range("E22").Value = "  -4.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.05"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.03"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.09"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.03%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("E27").Value = "  -0.86%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.84"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("E30").Value = "  -0.66%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.69"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.14"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.66%  "

$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.54"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.02%  "

$ws.Range("B34").Value = "Stacks"
$ws.Range("C34").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.32"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "55.47"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.34%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.91"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "455.09"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.87%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.202.08"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.66%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0799"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.71%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0385"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.57%  "

$ws.Range("E41").Value = "  -3.51%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.15"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.12%  "

$ws.Range("E43").Value = "  -11.87%  "

$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.42"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.08%  "

$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.246"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.04%  "

$ws.Range("E47").Value = "  -4.04%  "

$ws.Range("E48").Value = "  -1.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.35"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.79%  "

$ws.Range("E50").Value = "  -8.37%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.28"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.67%  "

